$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Name and Week ---
$ws.Range("B2").Value = "Jesse Hare"
$ws.Range("G2").Value = 10

# --- Activity column (A4:A8) filled first, in row order ---
$ws.Range("A4").Value = "front-end testing"
$ws.Range("A5").Value = "bug fixes"
$ws.Range("A6").Value = "client interview"
$ws.Range("A7").Value = "analyse new requirements"
$ws.Range("A8").Value = "implement client suggestions"

# --- Type column (C4:C8) - "G" rows filled before "I" rows so the
#     shared-string table ends up with G before I, matching the source file ---
$ws.Range("C6").Value = "G"
$ws.Range("C7").Value = "G"
$ws.Range("C8").Value = "G"
$ws.Range("C4").Value = "I"
$ws.Range("C5").Value = "I"

# --- Dates (D4:D8) ---
$ws.Range("D4").Value = "09/30/2019"
$ws.Range("D5").Value = "10/01/2019"
$ws.Range("D6").Value = "10/02/2019"
$ws.Range("D7").Value = "10/03/2019"
$ws.Range("D8").Value = "10/04/2019"

# --- Start times (E4:E8) - all 9:00 AM ---
$ws.Range("E4").Value = 0.375
$ws.Range("E5").Value = 0.375
$ws.Range("E6").Value = 0.375
$ws.Range("E7").Value = 0.375
$ws.Range("E8").Value = 0.375

# --- End times (F4:F8) ---
$ws.Range("F4").Value = 0.5
$ws.Range("F5").Value = 0.66666666666666663
$ws.Range("F6").Value = 0.45833333333333331
$ws.Range("F7").Value = 0.58333333333333337
$ws.Range("F8").Value = 1200

# --- Individual/Group hour totals ---
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 7
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 3

# --- Column widths: match the widths Excel settled on after the data entry
#     (ColumnWidth is in "characters"; the engine stores width = ColumnWidth + 5/6,
#     quantised to the nearest 1/6 character, so these inputs are chosen to land
#     on the closest representable value to the target stored widths of
#     22.42578125 / 13 / 13.7109375) ---
$ws.Columns("B").ColumnWidth = 21.6
$ws.Columns("D").ColumnWidth = 12.15
$ws.Columns("F").ColumnWidth = 12.8

# --- Move the active cell/selection ---
$ws.Range("D12").Select()
